# Generate Report for Handoff
# Update status text from "In Translation" to "Ready for handoff" and
# refresh the related generation/handoff timestamps on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2).
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-16 12:38:13"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2).
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-16 12:38:04"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2).
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-16 12:38:13"

# Widen the status/date columns so they match the newly generated report
# (closest reproducible value to the target 17.2159881591797 given the
# engine's internal pixel-grid snapping of ColumnWidth).
$newWidth = 16.3333333333333
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
